$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F3 and F4
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1084
$ws1.Range("F4").Value = 2501

# Sheet "全部类型" - update F5 and F6
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1084
$ws4.Range("F6").Value = 2501
